$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D/E values are set via Formula with a leading apostrophe (text-prefix)
# so numeric-looking strings (e.g. "1.07", "0.999") are stored as literal
# text instead of being coerced to numbers by Excels type inference.
# Style is reset to "Normal" right after so no stray formatting remains.
$textForced = @(
    @{ Cell = "D2"; Value = "63.696.29" }
    @{ Cell = "E2"; Value = "  -1.41%  " }
    @{ Cell = "D3"; Value = "3.079.66" }
    @{ Cell = "E3"; Value = "  -28.04%  " }
    @{ Cell = "D4"; Value = "0.999" }
    @{ Cell = "E4"; Value = "  -7.39%  " }
    @{ Cell = "D5"; Value = "592.73" }
    @{ Cell = "E5"; Value = "  +0.77%  " }
    @{ Cell = "D6"; Value = "155.37" }
    @{ Cell = "E6"; Value = "  +2.68%  " }
    @{ Cell = "E7"; Value = "  -8.52%  " }
    @{ Cell = "D8"; Value = "0.536" }
    @{ Cell = "E8"; Value = "  +1.60%  " }
    @{ Cell = "D9"; Value = "3.080.88" }
    @{ Cell = "E9"; Value = "  -2.19%  " }
    @{ Cell = "E10"; Value = "  -0.38%  " }
    @{ Cell = "D11"; Value = "5.93" }
    @{ Cell = "E11"; Value = "  +0.74%  " }
    @{ Cell = "D12"; Value = "0.451" }
    @{ Cell = "E12"; Value = "  -1.45%  " }
    @{ Cell = "D13"; Value = "0.0000238" }
    @{ Cell = "E13"; Value = "  -1.47%  " }
    @{ Cell = "D14"; Value = "36.66" }
    @{ Cell = "E14"; Value = "  -2.47%  " }
    @{ Cell = "E15"; Value = "  +1.01%  " }
    @{ Cell = "D16"; Value = "3.586.43" }
    @{ Cell = "E16"; Value = "  -1.36%  " }
    @{ Cell = "D17"; Value = "7.19" }
    @{ Cell = "E17"; Value = "  +0.14%  " }
    @{ Cell = "D18"; Value = "63.556.71" }
    @{ Cell = "E18"; Value = "  -2.65%  " }
    @{ Cell = "D19"; Value = "3.075.76" }
    @{ Cell = "E19"; Value = "  -1.57%  " }
    @{ Cell = "D20"; Value = "480.55" }
    @{ Cell = "E20"; Value = "  +2.67%  " }
    @{ Cell = "D21"; Value = "14.47" }
    @{ Cell = "E21"; Value = "  -2.37%  " }
    @{ Cell = "D22"; Value = "0.710" }
    @{ Cell = "E22"; Value = "  -3.23%  " }
    @{ Cell = "D23"; Value = "7.57" }
    @{ Cell = "E23"; Value = "  +0.45%  " }
    @{ Cell = "D24"; Value = "2.43" }
    @{ Cell = "E24"; Value = "  +2.62%  " }
    @{ Cell = "D25"; Value = "81.74" }
    @{ Cell = "E25"; Value = "  +0.30%  " }
    @{ Cell = "D26"; Value = "12.88" }
    @{ Cell = "E26"; Value = "  -2.21%  " }
    @{ Cell = "D27"; Value = "10.76" }
    @{ Cell = "E27"; Value = "  +10.28%  " }
    @{ Cell = "E28"; Value = "  -0.21%  " }
    @{ Cell = "D29"; Value = "7.66" }
    @{ Cell = "E29"; Value = "  +4.94%  " }
    @{ Cell = "D30"; Value = "2.69" }
    @{ Cell = "E30"; Value = "  -0.09%  " }
    @{ Cell = "D31"; Value = "0.999" }
    @{ Cell = "E31"; Value = "  -1.48%  " }
    @{ Cell = "D32"; Value = "2.20" }
    @{ Cell = "E32"; Value = "  -0.15%  " }
    @{ Cell = "E33"; Value = "  -3.20%  " }
    @{ Cell = "D34"; Value = "27.22" }
    @{ Cell = "E34"; Value = "  -0.19%  " }
    @{ Cell = "D35"; Value = "0.0₃0832" }
    @{ Cell = "E35"; Value = "  -2.16%  " }
    @{ Cell = "D36"; Value = "1.07" }
    @{ Cell = "E36"; Value = "  +0.96%  " }
    @{ Cell = "D37"; Value = "6.08" }
    @{ Cell = "E37"; Value = "  -0.44%  " }
    @{ Cell = "D38"; Value = "3.30" }
    @{ Cell = "E38"; Value = "  -1.49%  " }
    @{ Cell = "D39"; Value = "2.24" }
    @{ Cell = "E39"; Value = "  -0.97%  " }
    @{ Cell = "D40"; Value = "50.63" }
    @{ Cell = "E40"; Value = "  -0.65%  " }
    @{ Cell = "D41"; Value = "9.22" }
    @{ Cell = "E41"; Value = "  -0.50%  " }
    @{ Cell = "D42"; Value = "443.08" }
    @{ Cell = "E42"; Value = "  -1.36%  " }
    @{ Cell = "D43"; Value = "0.292" }
    @{ Cell = "E43"; Value = "  -0.40%  " }
    @{ Cell = "E44"; Value = "  +3.26%  " }
    @{ Cell = "D45"; Value = "0.0362" }
    @{ Cell = "E45"; Value = "  -2.04%  " }
    @{ Cell = "D46"; Value = "40.04" }
    @{ Cell = "E46"; Value = "  +2.73%  " }
    @{ Cell = "D47"; Value = "2.826.06" }
    @{ Cell = "E47"; Value = "  -1.26%  " }
    @{ Cell = "D48"; Value = "132.23" }
    @{ Cell = "E48"; Value = "  +1.02%  " }
    @{ Cell = "D49"; Value = "25.30" }
    @{ Cell = "E49"; Value = "  +0.74%  " }
    @{ Cell = "D50"; Value = "0.999" }
    @{ Cell = "E50"; Value = "  +0.01%  " }
    @{ Cell = "E51"; Value = "  -0.70%  " }
)

foreach ($u in $textForced) {
    $r = $ws.Range($u.Cell)
    $r.Formula = "'" + $u.Value
    $r.Style = "Normal"
}

# Plain text cells (coin names / URLs) - never numeric-like, safe to set directly.
$plainText = @(
    @{ Cell = "B35"; Value = "PEPE" }
    @{ Cell = "C35"; Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe" }
    @{ Cell = "B36"; Value = "Mantle" }
    @{ Cell = "C36"; Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt" }
    @{ Cell = "B45"; Value = "VeChain" }
    @{ Cell = "C45"; Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet" }
    @{ Cell = "B46"; Value = "Arweave" }
    @{ Cell = "C46"; Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar" }
    @{ Cell = "B49"; Value = "InjectiveProtocol" }
    @{ Cell = "C49"; Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj" }
    @{ Cell = "B50"; Value = "USDe" }
    @{ Cell = "C50"; Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde" }
)

foreach ($u in $plainText) {
    $ws.Range($u.Cell).Value = $u.Value
}
